$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Environments_OnGoing")

$ws.Range("A18").Value = "PartnersCommunityRelatedList"
$ws.Range("B18").Value = "/s/relatedlist/"

$ws.Range("B19").Value = "/AttachedContentDocuments"
$ws.Range("A19").Value = "PartnersCommunityRelatedListFiles"

$ws.Range("A20").Value = "PartnersCommunityRelatedListLineItems"
$ws.Range("B20").Value = "/OpportunityLineItems"

$ws.Range("B22").Select()
